# Color a set of "Things to do" list items green (RGB 00B050), matching
# the commit that highlighted the already-completed city/truck/menu work.
#
# RGB(0x00, 0xB0, 0x50) packed as a VBA-style BGR long == 0x50B000 == 5287936.
# Setting Range.Font.Color on a paragraph's Range (which includes the
# trailing paragraph mark) stamps <w:color> on the paragraph mark rPr
# (w:pPr/w:rPr) as well as on every run's rPr in that paragraph - exactly
# what the diff shows.

$d = $word.ActiveDocument
$green = 5287936

$targetTexts = @(
    "Create a Canvas for the main menu.",
    "Create button Start. ",
    "Create button Instructions.",
    "Create button Settings.",
    "Create a Game Manager as singleton.",
    "Truck.",
    "Street.",
    "Truck moving right.",
    "Truck moving left."
)

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text
    # Paragraph.Range.Text includes the trailing paragraph mark; trim it
    # (and any trailing cell-mark characters) before comparing.
    $trimmed = $text.TrimEnd([char]13, [char]7)
    if ($targetTexts -contains $trimmed) {
        $p.Range.Font.Color = $green
    }
}
